$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 123578.04
$ws.Range("I40").Value = 1000817
$ws.Range("J40").Value = 3954.5454
$ws.Range("K40").Value = 1000817
$ws.Range("L40").Value = 3954.5454
$ws.Range("M40").Value = -1000642
$ws.Range("N40").Value = -4304.5454

$ws.Range("H74").Value = 13998.6

$ws.Range("H77").Value = 13998.6

$ws.Range("H92").Value = 38462960
$ws.Range("I92").Value = 41668080
$ws.Range("J92").Value = 1505
$ws.Range("K92").Value = 41668080
$ws.Range("L92").Value = 1505
$ws.Range("M92").Value = -41666832
$ws.Range("N92").Value = -4001

$ws.Range("H98").Value = 1083
$ws.Range("I98").Value = 972.5833
$ws.Range("J98").Value = 1966.3334
$ws.Range("K98").Value = 972.5833
$ws.Range("L98").Value = 1966.3334
$ws.Range("M98").Value = 525.4167
$ws.Range("N98").Value = -4962.3334

$ws.Range("H122").Value = 1083
$ws.Range("I122").Value = 972.5833
$ws.Range("J122").Value = 1966.3334
$ws.Range("K122").Value = 2917.7499
$ws.Range("L122").Value = 5899.0002
$ws.Range("M122").Value = -467.7498999999998
$ws.Range("N122").Value = -10799.0002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 4455.778
$ws.Range("I45").Value = 4057.5
$ws.Range("K45").Value = 4057.5
$ws.Range("M45").Value = -3680.5

$ws.Range("H61").Value = 2297.9768
$ws.Range("I61").Value = 2114.7097
$ws.Range("K61").Value = 2114.7097
$ws.Range("M61").Value = -1902.7097

$ws.Range("H63").Value = 4744.8
$ws.Range("I63").Value = 3959.8
$ws.Range("K63").Value = 3959.8
$ws.Range("M63").Value = -3273.8

$ws.Range("H66").Value = 4744.8
$ws.Range("I66").Value = 3959.8
$ws.Range("K66").Value = 19799
$ws.Range("M66").Value = -16367

$ws.Range("H74").Value = 2829.5
$ws.Range("I74").Value = 2281.8484
$ws.Range("K74").Value = 2281.8484
$ws.Range("M74").Value = -1407.8484

$ws.Range("H77").Value = 2829.5
$ws.Range("I77").Value = 2281.8484
$ws.Range("K77").Value = 11409.242
$ws.Range("M77").Value = -7041.241999999998

$ws.Range("H97").Value = 882.2857
$ws.Range("I97").Value = 677.4545000000001
$ws.Range("K97").Value = 677.4545000000001
$ws.Range("M97").Value = -181.4545000000001

$ws.Range("H136").Value = 2297.9768
$ws.Range("I136").Value = 2114.7097
$ws.Range("K136").Value = 6344.1291
$ws.Range("M136").Value = -3794.1291

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2027
$ws.Range("I86").Value = 1697.6
$ws.Range("J86").Value = 2356.4
$ws.Range("K86").Value = 1697.6
$ws.Range("L86").Value = 2356.4
$ws.Range("M86").Value = -574.5999999999999
$ws.Range("N86").Value = -4602.4

$ws.Range("H89").Value = 2027
$ws.Range("I89").Value = 1697.6
$ws.Range("J89").Value = 2356.4
$ws.Range("K89").Value = 8488
$ws.Range("L89").Value = 11782
$ws.Range("M89").Value = -2872
$ws.Range("N89").Value = -23014

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2886.5293
$ws.Range("I31").Value = 2038.7273
$ws.Range("J31").Value = 4440.8335
$ws.Range("K31").Value = 2038.7273
$ws.Range("L31").Value = 4440.8335
$ws.Range("M31").Value = -1743.7273
$ws.Range("N31").Value = -5030.8335

$ws.Range("H34").Value = 2886.5293
$ws.Range("I34").Value = 2038.7273
$ws.Range("J34").Value = 4440.8335
$ws.Range("K34").Value = 2038.7273
$ws.Range("L34").Value = 4440.8335
$ws.Range("M34").Value = -1836.7273
$ws.Range("N34").Value = -4844.8335

$ws.Range("H58").Value = 2417.3784
$ws.Range("I58").Value = 2293.4062
$ws.Range("K58").Value = 2293.4062
$ws.Range("M58").Value = -2090.4062

$ws.Range("H107").Value = 85307.836
$ws.Range("I107").Value = 143884.86
$ws.Range("K107").Value = 143884.86
$ws.Range("M107").Value = -141964.86

$ws.Range("H134").Value = 2013.7
$ws.Range("I134").Value = 1686
$ws.Range("K134").Value = 5058
$ws.Range("M134").Value = -2523

$ws.Range("H136").Value = 2417.3784
$ws.Range("I136").Value = 2293.4062
$ws.Range("K136").Value = 6880.2186
$ws.Range("M136").Value = -4330.2186

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 107153896
$ws.Range("I4").Value = 193710100
$ws.Range("K4").Value = 581130300
$ws.Range("M4").Value = -581130188

$ws.Range("H122").Value = 498699.22
$ws.Range("J122").Value = 1162482.4
$ws.Range("L122").Value = 10462341.6
$ws.Range("N122").Value = -10467241.6

$ws.Range("H129").Value = 2222.5
$ws.Range("I129").Value = 844.5
$ws.Range("J129").Value = 3141.1667
$ws.Range("K129").Value = 2533.5
$ws.Range("L129").Value = 9423.500100000001
$ws.Range("M129").Value = 2466.5
$ws.Range("N129").Value = -19423.5001

$ws.Range("H140").Value = 1761.7646
$ws.Range("I140").Value = 1460.7142
$ws.Range("J140").Value = 3166.6667
$ws.Range("K140").Value = 4382.142599999999
$ws.Range("L140").Value = 9500.000100000001
$ws.Range("M140").Value = 797.8574000000008
$ws.Range("N140").Value = -19860.0001

$ws.Range("H141").Value = 3529.6365
$ws.Range("I141").Value = 3529.6365
$ws.Range("K141").Value = 10588.9095
$ws.Range("M141").Value = -5408.9095

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H109").Value = 100000
$ws.Range("J109").Value = 100000
$ws.Range("L109").Value = 100000
$ws.Range("N109").Value = -102080

$ws.Range("H122").Value = 3564.0833
$ws.Range("I122").Value = 4431.6665
$ws.Range("K122").Value = 13294.9995
$ws.Range("M122").Value = -10844.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 6695.4517
$ws.Range("I46").Value = 8433.625
$ws.Range("K46").Value = 8433.625
$ws.Range("M46").Value = -8245.625

$ws.Range("H93").Value = 37038100
$ws.Range("I93").Value = 62500660
$ws.Range("K93").Value = 62500660
$ws.Range("M93").Value = -62499412

$ws.Range("H123").Value = 29888
$ws.Range("J123").Value = 29888
$ws.Range("L123").Value = 29888
$ws.Range("N123").Value = -39688

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 666
$ws.Range("I107").Value = 499
$ws.Range("K107").Value = 1497
$ws.Range("M107").Value = 423

$ws.Range("H122").Value = 62503588
$ws.Range("I122").Value = 111111576
$ws.Range("K122").Value = 333334728
$ws.Range("M122").Value = -333332278

$ws.Range("H126").Value = 5907.3335
$ws.Range("I126").Value = 6298.8
$ws.Range("K126").Value = 18896.4
$ws.Range("M126").Value = -16426.4
